# Update the NATMI LR-pair output (Ifnk-Ifnar1) with refreshed TPM-derived
# values. This recomputation changes several numeric columns (G,H,I,J,
# M,N,O,P,Q,R,S,T) and, for the "MuSCs" sending-cluster rows, the
# ligand-expressing-cell count/rate (E,F) as well, per the new TPM input.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> ECs ---
$ws.Range("G2").Value = 0.062802
$ws.Range("H2").Value = 0.188406
$ws.Range("I2").Value = 0.0482629354765083
$ws.Range("J2").Value = 0.0482629354765083
$ws.Range("M2").Value = 10.70375566666667
$ws.Range("N2").Value = 32.111267
$ws.Range("O2").Value = 0.3178747217938744
$ws.Range("P2").Value = 0.3178747217938744
$ws.Range("Q2").Value = 0.6722172633779999
$ws.Range("R2").Value = 6.049955370401999
$ws.Range("S2").Value = 0.01534156718755079
$ws.Range("T2").Value = 0.01534156718755079

# --- Row 3: ECs -> FAPs ---
$ws.Range("G3").Value = 0.062802
$ws.Range("H3").Value = 0.188406
$ws.Range("I3").Value = 0.0482629354765083
$ws.Range("J3").Value = 0.0482629354765083
$ws.Range("O3").Value = 0.4731463873433873
$ws.Range("P3").Value = 0.4731463873433873
$ws.Range("Q3").Value = 1.000573961598
$ws.Range("R3").Value = 9.005165654381999
$ws.Range("S3").Value = 0.0228354335632969
$ws.Range("T3").Value = 0.0228354335632969

# --- Row 4: ECs -> MuSCs ---
$ws.Range("G4").Value = 0.062802
$ws.Range("H4").Value = 0.188406
$ws.Range("I4").Value = 0.0482629354765083
$ws.Range("J4").Value = 0.0482629354765083
$ws.Range("M4").Value = 7.036919999999999
$ws.Range("N4").Value = 21.11076
$ws.Range("O4").Value = 0.2089788908627383
$ws.Range("P4").Value = 0.2089788908627384
$ws.Range("Q4").Value = 0.4419326498399999
$ws.Range("R4").Value = 3.97739384856
$ws.Range("S4").Value = 0.01008593472566061
$ws.Range("T4").Value = 0.01008593472566061

# --- Row 5: FAPs -> ECs ---
$ws.Range("G5").Value = 0.5410386666666667
$ws.Range("I5").Value = 0.415784756212054
$ws.Range("J5").Value = 0.415784756212054
$ws.Range("M5").Value = 10.70375566666667
$ws.Range("N5").Value = 32.111267
$ws.Range("O5").Value = 0.3178747217938744
$ws.Range("P5").Value = 0.3178747217938744
$ws.Range("Q5").Value = 5.791145694219111
$ws.Range("R5").Value = 52.120311247972
$ws.Range("S5").Value = 0.1321674637070406
$ws.Range("T5").Value = 0.1321674637070406

# --- Row 6: FAPs -> FAPs ---
$ws.Range("G6").Value = 0.5410386666666667
$ws.Range("I6").Value = 0.415784756212054
$ws.Range("J6").Value = 0.415784756212054
$ws.Range("O6").Value = 0.4731463873433873
$ws.Range("P6").Value = 0.4731463873433873
$ws.Range("Q6").Value = 8.619935704028
$ws.Range("R6").Value = 77.579421336252
$ws.Range("S6").Value = 0.1967270553141844
$ws.Range("T6").Value = 0.1967270553141844

# --- Row 7: FAPs -> MuSCs ---
$ws.Range("G7").Value = 0.5410386666666667
$ws.Range("I7").Value = 0.415784756212054
$ws.Range("J7").Value = 0.415784756212054
$ws.Range("M7").Value = 7.036919999999999
$ws.Range("N7").Value = 21.11076
$ws.Range("O7").Value = 0.2089788908627383
$ws.Range("P7").Value = 0.2089788908627384
$ws.Range("Q7").Value = 3.80724581424
$ws.Range("R7").Value = 34.26521232816
$ws.Range("S7").Value = 0.0868902371908291
$ws.Range("T7").Value = 0.08689023719082911

# --- Row 8: MuSCs -> ECs ---
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.6974063333333334
$ws.Range("H8").Value = 2.092219
$ws.Range("I8").Value = 0.5359523083114377
$ws.Range("J8").Value = 0.5359523083114377
$ws.Range("M8").Value = 10.70375566666667
$ws.Range("N8").Value = 32.111267
$ws.Range("O8").Value = 0.3178747217938744
$ws.Range("P8").Value = 0.3178747217938744
$ws.Range("Q8").Value = 7.464866992385889
$ws.Range("R8").Value = 67.183802931473
$ws.Range("S8").Value = 0.1703656908992831
$ws.Range("T8").Value = 0.1703656908992831

# --- Row 9: MuSCs -> FAPs ---
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.6974063333333334
$ws.Range("H9").Value = 2.092219
$ws.Range("I9").Value = 0.5359523083114377
$ws.Range("J9").Value = 0.5359523083114377
$ws.Range("O9").Value = 0.4731463873433873
$ws.Range("P9").Value = 0.4731463873433873
$ws.Range("Q9").Value = 11.111216486527
$ws.Range("R9").Value = 100.000948378743
$ws.Range("S9").Value = 0.253583898465906
$ws.Range("T9").Value = 0.2535838984659061

# --- Row 10: MuSCs -> MuSCs ---
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.6974063333333334
$ws.Range("H10").Value = 2.092219
$ws.Range("I10").Value = 0.5359523083114377
$ws.Range("J10").Value = 0.5359523083114377
$ws.Range("M10").Value = 7.036919999999999
$ws.Range("N10").Value = 21.11076
$ws.Range("O10").Value = 0.2089788908627383
$ws.Range("P10").Value = 0.2089788908627384
$ws.Range("Q10").Value = 4.90759257516
$ws.Range("R10").Value = 44.16833317644
$ws.Range("S10").Value = 0.1120027189462486
$ws.Range("T10").Value = 0.1120027189462486
